$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (header "K") values for rows 2-8 with recalculated strikeout counts
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 1
